$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for columns B..F, rows 2..15 (row index -> values)
$data = @{
    2  = @("NSE:APCOTEXIND", "NSE:BANDHANBNK", "NSE:HINDPETRO", "NSE:BEL", "NSE:ASIANPAINT")
    3  = @("NSE:BLKASHYAP", "NSE:DBCORP", "", "NSE:CROMPTON", "NSE:PAGEIND")
    4  = @("NSE:IPL", "NSE:FDC", "", "NSE:NAUKRI", "")
    5  = @("NSE:MOL", "NSE:IMAGICAA", "", "NSE:PETRONET", "")
    6  = @("NSE:PAGEIND", "NSE:INTENTECH", "", "", "")
    7  = @("NSE:PGHL", "NSE:JISLJALEQS", "", "", "")
    8  = @("", "NSE:NGIL", "", "", "")
    9  = @("", "NSE:NH", "", "", "")
    10 = @("", "NSE:PHOENIXLTD", "", "", "")
    11 = @("", "NSE:PPL", "", "", "")
    12 = @("", "NSE:PTL", "", "", "")
    13 = @("", "NSE:RELCHEMQ", "", "", "")
    14 = @("", "NSE:RITES", "", "", "")
    15 = @("", "NSE:ROSSARI", "", "", "")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}

# Remove now-unused rows 16..27 entirely
$ws.Rows("16:27").Delete()
